# Actualización desde MV -datos-
# Adds 5 new daily rows (20-09-2021 .. 24-09-2021) to the EMBI spreads sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("20-09-2021", 322.4, 395, 205, 317, 85.2, 120.9, 12, 154.8, 498.1, 1575.9, 298.9, 291, 150, 357, 177),
    @("21-09-2021", 320.6, 393, 205, 312, 85.2, 119.8, 16.1, 151.9, 487.7, 1596.5, 295.9, 288, 147, 352, 175),
    @("22-09-2021", 321.8, 394, 205, 314, 84.6, 118.9, 14.7, 152, 492.7, 1612.3, 295.7, 290, 146, 354, 175),
    @("23-09-2021", 317.3, 390, 198, 313, 84, 115, 11.6, 147.1, 501.2, 1598.2, 294.3, 286, 141, 351, 169),
    @("24-09-2021", 319.9, 394, 199, 318, 85.3, 119.1, 12.8, 148.3, 510.4, 1619, 297.3, 289, 141, 355, 170)
)

$startRow = 181
for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowValues = $newRows[$r]
    $rowIndex = $startRow + $r
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $rowValues[$c]
    }
}
